$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna4"
$ws.Cells.Item(2, 3).Value = "Epha3"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.6731353333333333
$ws.Cells.Item(2, 8).Value = 2.019406
$ws.Cells.Item(2, 9).Value = 0.3272865828458516
$ws.Cells.Item(2, 10).Value = 0.3272865828458516
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 46.25093466666667
$ws.Cells.Item(2, 14).Value = 138.752804
$ws.Cells.Item(2, 15).Value = 0.9569015955251317
$ws.Cells.Item(2, 16).Value = 0.9569015955251318
$ws.Cells.Item(2, 17).Value = 31.13313832382489
$ws.Cells.Item(2, 18).Value = 280.198244914424
$ws.Cells.Item(2, 19).Value = 0.3131810533191636
$ws.Cells.Item(2, 20).Value = 0.3131810533191636
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna4"
$ws.Cells.Item(3, 3).Value = "Epha3"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.6731353333333333
$ws.Cells.Item(3, 8).Value = 2.019406
$ws.Cells.Item(3, 9).Value = 0.3272865828458516
$ws.Cells.Item(3, 10).Value = 0.3272865828458516
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.01632333333333333
$ws.Cells.Item(3, 14).Value = 0.04897
$ws.Cells.Item(3, 15).Value = 0.0003377190931065126
$ws.Cells.Item(3, 16).Value = 0.0003377190931065127
$ws.Cells.Item(3, 17).Value = 0.01098781242444444
$ws.Cells.Item(3, 18).Value = 0.09889031182000001
$ws.Cells.Item(3, 19).Value = 0.0001105309279446305
$ws.Cells.Item(3, 20).Value = 0.0001105309279446305
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna4"
$ws.Cells.Item(4, 3).Value = "Epha3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.6731353333333333
$ws.Cells.Item(4, 8).Value = 2.019406
$ws.Cells.Item(4, 9).Value = 0.3272865828458516
$ws.Cells.Item(4, 10).Value = 0.3272865828458516
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.066797333333333
$ws.Cells.Item(4, 14).Value = 6.200391999999999
$ws.Cells.Item(4, 15).Value = 0.04276068538176181
$ws.Cells.Item(4, 16).Value = 0.04276068538176181
$ws.Cells.Item(4, 17).Value = 1.391234311905778
$ws.Cells.Item(4, 18).Value = 12.521108807152
$ws.Cells.Item(4, 19).Value = 0.01399499859874338
$ws.Cells.Item(4, 20).Value = 0.01399499859874338
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efna4"
$ws.Cells.Item(5, 3).Value = "Epha3"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.9964423333333334
$ws.Cells.Item(5, 8).Value = 2.989327
$ws.Cells.Item(5, 9).Value = 0.4844823769162027
$ws.Cells.Item(5, 10).Value = 0.4844823769162026
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 46.25093466666667
$ws.Cells.Item(5, 14).Value = 138.752804
$ws.Cells.Item(5, 15).Value = 0.9569015955251317
$ws.Cells.Item(5, 16).Value = 0.9569015955251318
$ws.Cells.Item(5, 17).Value = 46.08638925810089
$ws.Cells.Item(5, 18).Value = 414.777503322908
$ws.Cells.Item(5, 19).Value = 0.4636019594749226
$ws.Cells.Item(5, 20).Value = 0.4636019594749226
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efna4"
$ws.Cells.Item(6, 3).Value = "Epha3"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.9964423333333334
$ws.Cells.Item(6, 8).Value = 2.989327
$ws.Cells.Item(6, 9).Value = 0.4844823769162027
$ws.Cells.Item(6, 10).Value = 0.4844823769162026
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.01632333333333333
$ws.Cells.Item(6, 14).Value = 0.04897
$ws.Cells.Item(6, 15).Value = 0.0003377190931065126
$ws.Cells.Item(6, 16).Value = 0.0003377190931065127
$ws.Cells.Item(6, 17).Value = 0.01626526035444445
$ws.Cells.Item(6, 18).Value = 0.14638734319
$ws.Cells.Item(6, 19).Value = 0.0001636189489582276
$ws.Cells.Item(6, 20).Value = 0.0001636189489582276
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna4"
$ws.Cells.Item(7, 3).Value = "Epha3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.9964423333333334
$ws.Cells.Item(7, 8).Value = 2.989327
$ws.Cells.Item(7, 9).Value = 0.4844823769162027
$ws.Cells.Item(7, 10).Value = 0.4844823769162026
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.066797333333333
$ws.Cells.Item(7, 14).Value = 6.200391999999999
$ws.Cells.Item(7, 15).Value = 0.04276068538176181
$ws.Cells.Item(7, 16).Value = 0.04276068538176181
$ws.Cells.Item(7, 17).Value = 2.059444357353778
$ws.Cells.Item(7, 18).Value = 18.534999216184
$ws.Cells.Item(7, 19).Value = 0.02071679849232188
$ws.Cells.Item(7, 20).Value = 0.02071679849232188
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Efna4"
$ws.Cells.Item(8, 3).Value = "Epha3"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.3871376666666667
$ws.Cells.Item(8, 8).Value = 1.161413
$ws.Cells.Item(8, 9).Value = 0.1882310402379457
$ws.Cells.Item(8, 10).Value = 0.1882310402379457
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 46.25093466666667
$ws.Cells.Item(8, 14).Value = 138.752804
$ws.Cells.Item(8, 15).Value = 0.9569015955251317
$ws.Cells.Item(8, 16).Value = 0.9569015955251318
$ws.Cells.Item(8, 17).Value = 17.90547892800578
$ws.Cells.Item(8, 18).Value = 161.149310352052
$ws.Cells.Item(8, 19).Value = 0.1801185827310455
$ws.Cells.Item(8, 20).Value = 0.1801185827310456
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Efna4"
$ws.Cells.Item(9, 3).Value = "Epha3"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.3871376666666667
$ws.Cells.Item(9, 8).Value = 1.161413
$ws.Cells.Item(9, 9).Value = 0.1882310402379457
$ws.Cells.Item(9, 10).Value = 0.1882310402379457
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.01632333333333333
$ws.Cells.Item(9, 14).Value = 0.04897
$ws.Cells.Item(9, 15).Value = 0.0003377190931065126
$ws.Cells.Item(9, 16).Value = 0.0003377190931065127
$ws.Cells.Item(9, 17).Value = 0.006319377178888888
$ws.Cells.Item(9, 18).Value = 0.05687439461
$ws.Cells.Item(9, 19).Value = 0.00006356921620365452
$ws.Cells.Item(9, 20).Value = 0.00006356921620365454
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Efna4"
$ws.Cells.Item(10, 3).Value = "Epha3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.3871376666666667
$ws.Cells.Item(10, 8).Value = 1.161413
$ws.Cells.Item(10, 9).Value = 0.1882310402379457
$ws.Cells.Item(10, 10).Value = 0.1882310402379457
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.066797333333333
$ws.Cells.Item(10, 14).Value = 6.200391999999999
$ws.Cells.Item(10, 15).Value = 0.04276068538176181
$ws.Cells.Item(10, 16).Value = 0.04276068538176181
$ws.Cells.Item(10, 17).Value = 0.8001350970995554
$ws.Cells.Item(10, 18).Value = 7.201215873895999
$ws.Cells.Item(10, 19).Value = 0.008048888290696546
$ws.Cells.Item(10, 20).Value = 0.008048888290696546